$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.181499361991882
$ws.Range("B1").Value = 2.257375717163086
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.349781036376953
$ws.Range("E1").Value = 1.22279167175293
